$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.787.71"
$ws.Range("E2").Value = "  -0.48%  "
$ws.Range("D3").Value = "3.988.75"
$ws.Range("E3").Value = "  -1.20%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "540.91"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +4.09%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "149.39"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.05%  "
$ws.Range("E7").Value = "  +12.09%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.00"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.741"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.93%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.169"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.72%  "
$ws.Range("E11").Value = "  -3.37%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "49.84"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +4.36%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.62"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -1.81%  "
$ws.Range("D14").Value = "4.628.47"
$ws.Range("E14").Value = "  -1.09%  "
$ws.Range("D15").Value = "3.992.18"
$ws.Range("E15").Value = "  -1.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.99"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -0.67%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "20.35"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -3.81%  "
$ws.Range("E18").Value = "  -0.22%  "
$ws.Range("E19").Value = "  -2.53%  "
$ws.Range("D20").Value = "71.759.05"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "426.87"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "96.85"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.50%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.48"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -1.32%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "4.23"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +5.52%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "14.20"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -3.13%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "11.32"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -5.19%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.61"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -5.12%  "
$ws.Range("B28").Value = "LEO"
$ws.Range("C28").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "5.85"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +1.29%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "36.62"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  -0.90%  "
$ws.Range("B30").Value = "Toncoin"
$ws.Range("C30").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.62"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +16.67%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.130"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +1.21%  "
$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "13.29"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -1.62%  "
$ws.Range("B33").Value = "NEARProtocol"
$ws.Range("C33").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.17"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.08%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "48.43"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +19.53%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "672.98"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -3.23%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "65.43"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -4.57%  "
$ws.Range("E37").Value = "  +0.67%  "
$ws.Range("E38").Value = "  -2.47%  "
$ws.Range("E39").Value = "  -8.50%  "
$ws.Range("E40").Value = "  -6.93%  "
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.33"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +5.20%  "
$ws.Range("E43").Value = "  +0.36%  "
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("E45").Value = "  +2.62%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.79"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +8.33%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.66"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -3.44%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.35"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -4.67%  "
$ws.Range("E49").Value = "  -3.90%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000273"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "144.57"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.24%  "
